$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.347.49'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '3.398.83'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '179.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.199'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +8.00%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '48.31'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('E12').Value = '  +3.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '681.45'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').Value = '3.947.45'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.60'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '69.469.82'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').Value = '3.397.80'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.69'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.28'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.08'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '101.12'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.74'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '33.50'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.76'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.91'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.78'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +13.06%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '555.32'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.106'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '58.05'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '3.608.73'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.141'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.30'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.25%  '
$ws.Range('D40').Value = '0.0₃0748'
$ws.Range('E40').Value = '  +11.06%  '
$ws.Range('E41').Value = '  +4.45%  '
$ws.Range('E42').Value = '  +3.31%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0427'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.70%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.336'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.68'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.129'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.39'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.04%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '131.01'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('B50').Value = 'CoreDAO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.62'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.43'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.30%  '
